$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = "4/21/2021"
